# Weekly update: insert two new observation rows (228-229) for
# Fruta / Feria Lagunitas de Puerto Montt - Pera, pushing the existing
# rows 228-314 down to 230-316.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 228; Excel's normal
# row-insert behaviour copies formatting (incl. the date number format
# in column D) from the row above, same as interactively inserting rows.
$ws.Rows(228).Insert()
$ws.Rows(228).Insert()

# --- New row 228: Packham's Triumph, Primera ---
$ws.Range("A228").Value = 4
$ws.Range("B228").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C228").Value = "Los Lagos"
$ws.Range("D228").Value = 44784
$ws.Range("E228").Value = 10
$ws.Range("F228").Value = "Fruta"
$ws.Range("G228").Value = 100104
$ws.Range("H228").Value = "Frutos de pepita"
$ws.Range("I228").Value = 100104005
$ws.Range("J228").Value = "Pera"
$ws.Range("K228").Value = "Packham's Triumph"
$ws.Range("L228").Value = "Primera"
$ws.Range("M228").Value = 300
$ws.Range("N228").Value = 15000
$ws.Range("O228").Value = 16000
$ws.Range("P228").Value = 15500
$ws.Range("Q228").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R228").Value = "Región de O'Higgins"
$ws.Range("S228").Value = 1033
$ws.Range("T228").Value = 15

# --- New row 229: Packham's Triumph, Segunda ---
$ws.Range("A229").Value = 4
$ws.Range("B229").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value = "Los Lagos"
$ws.Range("D229").Value = 44784
$ws.Range("E229").Value = 10
$ws.Range("F229").Value = "Fruta"
$ws.Range("G229").Value = 100104
$ws.Range("H229").Value = "Frutos de pepita"
$ws.Range("I229").Value = 100104005
$ws.Range("J229").Value = "Pera"
$ws.Range("K229").Value = "Packham's Triumph"
$ws.Range("L229").Value = "Segunda"
$ws.Range("M229").Value = 150
$ws.Range("N229").Value = 13000
$ws.Range("O229").Value = 13000
$ws.Range("P229").Value = 13000
$ws.Range("Q229").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R229").Value = "Región de O'Higgins"
$ws.Range("S229").Value = 867
$ws.Range("T229").Value = 15
